$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.093.61"
$ws.Range("E2").Value = "  -3.60%  "

# Row 3
$ws.Range("D3").Value = "3.680.21"
$ws.Range("E3").Value = "  -2.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "590.96"
$ws.Range("E5").Value = "  -3.91%  "

# Row 6
$ws.Range("D6").Value = "166.86"
$ws.Range("E6").Value = "  -6.11%  "

# Row 7
$ws.Range("D7").Value = "3.680.13"
$ws.Range("E7").Value = "  -2.95%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -1.05%  "

# Row 10
$ws.Range("E10").Value = "  -3.35%  "

# Row 11
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").Value = "  -5.28%  "

# Row 12
$ws.Range("E12").Value = "  -5.65%  "

# Row 13
$ws.Range("D13").Value = "37.67"
$ws.Range("E13").Value = "  -5.56%  "

# Row 14
$ws.Range("D14").Value = "0.0000239"
$ws.Range("E14").Value = "  -6.07%  "

# Row 15
$ws.Range("D15").Value = "4.296.12"
$ws.Range("E15").Value = "  -2.64%  "

# Row 16
$ws.Range("D16").Value = "3.680.55"
$ws.Range("E16").Value = "  -2.57%  "

# Row 17
$ws.Range("D17").Value = "67.148.53"
$ws.Range("E17").Value = "  -3.61%  "

# Row 18
$ws.Range("E18").Value = "  -4.17%  "

# Row 19
$ws.Range("E19").Value = "  -6.50%  "

# Row 20
$ws.Range("D20").Value = "16.98"
$ws.Range("E20").Value = "  +2.02%  "

# Row 21
$ws.Range("D21").Value = "485.77"
$ws.Range("E21").Value = "  -4.93%  "

# Row 22
$ws.Range("D22").Value = "9.10"
$ws.Range("E22").Value = "  -4.76%  "

# Row 23
$ws.Range("D23").Value = "0.718"
$ws.Range("E23").Value = "  -2.54%  "

# Row 24
$ws.Range("D24").Value = "84.76"
$ws.Range("E24").Value = "  -1.92%  "

# Row 25
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  -7.74%  "

# Row 26
$ws.Range("E26").Value = "  -1.26%  "

# Row 27
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").Value = "  -6.56%  "

# Row 28
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.53%  "

# Row 29
$ws.Range("D29").Value = "9.92"
$ws.Range("E29").Value = "  -6.22%  "

# Row 30
$ws.Range("E30").Value = "  -3.49%  "

# Row 31
$ws.Range("E31").Value = "  -7.04%  "

# Row 32
$ws.Range("D32").Value = "7.68"
$ws.Range("E32").Value = "  -5.61%  "

# Row 33
$ws.Range("D33").Value = "31.72"
$ws.Range("E33").Value = "  +1.27%  "

# Row 34
$ws.Range("D34").Value = "3.819.87"
$ws.Range("E34").Value = "  -2.63%  "

# Row 35
$ws.Range("D35").Value = "3.617.58"
$ws.Range("E35").Value = "  -2.59%  "

# Row 36
$ws.Range("D36").Value = "0.106"
$ws.Range("E36").Value = "  -7.57%  "

# Row 38
$ws.Range("E38").Value = "  -5.37%  "

# Row 39
$ws.Range("E39").Value = "  -7.16%  "

# Row 40
$ws.Range("D40").Value = "0.130"
$ws.Range("E40").Value = "  -8.12%  "

# Row 41
$ws.Range("E41").Value = "  -6.20%  "

# Row 42
$ws.Range("D42").Value = "440.60"
$ws.Range("E42").Value = "  -7.69%  "

# Row 43
$ws.Range("D43").Value = "48.69"
$ws.Range("E43").Value = "  -2.21%  "

# Row 44
$ws.Range("D44").Value = "1.94"
$ws.Range("E44").Value = "  -6.23%  "

# Row 45
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  -8.30%  "

# Row 46
$ws.Range("D46").Value = "8.27"
$ws.Range("E46").Value = "  -3.90%  "

# Row 47
$ws.Range("E47").Value = "  +0.00%  "

# Row 48
$ws.Range("D48").Value = "39.74"
$ws.Range("E48").Value = "  -10.33%  "

# Row 49
$ws.Range("D49").Value = "140.65"

# Row 50
$ws.Range("D50").Value = "2.764.01"
$ws.Range("E50").Value = "  -6.30%  "

# Row 51
$ws.Range("E51").Value = "  -5.19%  "
